$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap data between row 2 and row 3 for columns A, B, D, E, F, G, H, Q, R
# and move the "Hack" value in column AI from row 3 to row 2.

$ws.Range("A2").Value = 89819259
$ws.Range("B2").Value = 56411
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 100049
$ws.Range("F2").Value = "Spillkråka"
$ws.Range("G2").Value = "Dryocopus martius"
$ws.Range("H2").Value = "(Linnaeus, 1758)"
$ws.Range("Q2").Value = 584738.7995661208
$ws.Range("R2").Value = 6695804.92004218
$ws.Range("AI2").Value = "Hack"

$ws.Range("A3").Value = 89819263
$ws.Range("B3").Value = 43464
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 101735
$ws.Range("F3").Value = "Jättesvampmal"
$ws.Range("G3").Value = "Scardia boletella"
$ws.Range("H3").Value = "(Fabricius, 1794)"
$ws.Range("Q3").Value = 584567.934929442
$ws.Range("R3").Value = 6696037.803599558
$ws.Range("AI3").ClearContents()
